# ============================================================================
# Update gh-pages output (广州-漫展信息.xlsx) to the data generated at a3196b5
# ----------------------------------------------------------------------------
# Sheets:
#   1 展览     (exhibitions)
#   2 演出     (live shows)
#   3 本地生活 (local life - empty, untouched)
#   4 全部类型 (all types - static merged snapshot of sheets 1+2+3 sorted by date)
# ============================================================================

$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------------
# Sheet "展览" (sheet1): refresh "想去人数" (want-to-go count) figures, and
# mark the 5th Qingyun Comic Con as cancelled.
# ----------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")

$ws1.Cells.Item(2,6).Value = 11369
$ws1.Cells.Item(3,6).Value = 1871
$ws1.Cells.Item(4,6).Value = 519
$ws1.Cells.Item(5,6).Value = 816
$ws1.Cells.Item(6,6).Value = 2411
$ws1.Cells.Item(7,6).Value = 737
$ws1.Cells.Item(8,6).Value = 986
$ws1.Cells.Item(9,6).Value = 577
$ws1.Cells.Item(10,6).Value = 445
$ws1.Cells.Item(11,6).Value = 475
$ws1.Cells.Item(12,6).Value = 437
$ws1.Cells.Item(13,6).Value = 1299
$ws1.Cells.Item(14,6).Value = 628
$ws1.Cells.Item(15,6).Value = 80
$ws1.Cells.Item(16,6).Value = 957
$ws1.Cells.Item(17,6).Value = 482
$ws1.Cells.Item(18,6).Value = 638
$ws1.Cells.Item(19,6).Value = 1031
$ws1.Cells.Item(20,6).Value = 196
$ws1.Cells.Item(21,6).Value = 913
$ws1.Cells.Item(22,6).Value = 121
$ws1.Cells.Item(23,6).Value = 237
$ws1.Cells.Item(24,6).Value = 113
$ws1.Cells.Item(25,6).Value = 250
$ws1.Cells.Item(26,6).Value = 657
$ws1.Cells.Item(27,6).Value = 156
$ws1.Cells.Item(28,6).Value = 97
$ws1.Cells.Item(29,6).Value = 312

# Row 24 - 广州·第五届清云动漫展 got cancelled
$ws1.Cells.Item(24,3).Value = "广州·第五届清云动漫展（取消）"
$ws1.Cells.Item(24,7).Value = "不可售"

# ----------------------------------------------------------------------------
# Sheet "演出" (sheet2): bump 次元LAB's want-to-go count, insert the newly
# announced HANAPOKO show (2024-03-09) ahead of the "三月的幻想" show, and
# bump "春卷饭十周年" want-to-go count.
# ----------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")

# 广州·次元LAB 二次元电音节 : 848 -> 853
$ws2.Cells.Item(5,6).Value = 853

# Insert a new row 7 (shifts 三月的幻想/春卷饭/夏川里美 down by one row each)
$ws2.Rows.Item(7).Insert()

# Copy column-A number formatting/style from the row above onto the new row
$ws2.Cells.Item(6,1).Copy()
$ws2.Cells.Item(7,1).PasteSpecial(-4122)  # xlPasteFormats

$ws2.Cells.Item(7,1).Value = 6
$ws2.Cells.Item(7,2).Value = "'2024-03-09"
$ws2.Cells.Item(7,3).Value = "【大会员抢先购】广州·HANAPOKO 2024 LIVE"
$ws2.Cells.Item(7,4).Value = "海珠同创汇东一街11号（上冲南约11-2） 声音共和Livehouse"
$ws2.Cells.Item(7,5).Value = "2024.03.09 14:00-03.09 15:30"
$ws2.Cells.Item(7,6).Value = 1
$ws2.Cells.Item(7,7).Value = "'380"
$ws2.Cells.Item(7,8).Value = $false
$ws2.Cells.Item(7,9).Value = "https://show.bilibili.com/platform/detail.html?id=81279"
$ws2.Cells.Item(7,10).Value = "//i2.hdslb.com/bfs/openplatform/202401/tMZ1Jp2G1705992352054.jpeg"

# 广州·春卷饭 十周年 2024 专场演出 (now on row 9 after the insert above) : 248 -> 300
$ws2.Cells.Item(9,6).Value = 300

# ----------------------------------------------------------------------------
# Sheet "全部类型" (sheet4): static merged snapshot of 展览+演出+本地生活 sorted
# by date - mirror the same edits made above in their merged-sheet position.
# ----------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Cells.Item(2,6).Value = 11369
$ws4.Cells.Item(3,6).Value = 1871
$ws4.Cells.Item(5,6).Value = 519
$ws4.Cells.Item(6,6).Value = 816
$ws4.Cells.Item(7,6).Value = 2411
$ws4.Cells.Item(8,6).Value = 737
$ws4.Cells.Item(9,6).Value = 986
$ws4.Cells.Item(11,6).Value = 577
$ws4.Cells.Item(12,6).Value = 445
$ws4.Cells.Item(13,6).Value = 475
$ws4.Cells.Item(14,6).Value = 437
$ws4.Cells.Item(15,6).Value = 1299
$ws4.Cells.Item(17,6).Value = 628
$ws4.Cells.Item(18,6).Value = 80
$ws4.Cells.Item(19,6).Value = 853
$ws4.Cells.Item(20,6).Value = 957
$ws4.Cells.Item(21,6).Value = 482
$ws4.Cells.Item(22,6).Value = 638
$ws4.Cells.Item(23,6).Value = 1031
$ws4.Cells.Item(24,6).Value = 196
$ws4.Cells.Item(25,6).Value = 913
$ws4.Cells.Item(26,6).Value = 121
$ws4.Cells.Item(27,6).Value = 237
$ws4.Cells.Item(29,6).Value = 113
$ws4.Cells.Item(30,6).Value = 250
$ws4.Cells.Item(31,6).Value = 657
$ws4.Cells.Item(32,6).Value = 156
$ws4.Cells.Item(34,6).Value = 97
$ws4.Cells.Item(36,6).Value = 312

# Row 29 - 广州·第五届清云动漫展 got cancelled
$ws4.Cells.Item(29,3).Value = "广州·第五届清云动漫展（取消）"
$ws4.Cells.Item(29,7).Value = "不可售"

# Insert the new HANAPOKO row at its sorted-by-date position (row 31, between
# 明日方舟ONLY on 03-02 and 代号鸢only2.0 on 03-16)
$ws4.Rows.Item(31).Insert()

$ws4.Cells.Item(30,1).Copy()
$ws4.Cells.Item(31,1).PasteSpecial(-4122)  # xlPasteFormats

$ws4.Cells.Item(31,1).Value = 30
$ws4.Cells.Item(31,2).Value = "'2024-03-09"
$ws4.Cells.Item(31,3).Value = "【大会员抢先购】广州·HANAPOKO 2024 LIVE"
$ws4.Cells.Item(31,4).Value = "海珠同创汇东一街11号（上冲南约11-2） 声音共和Livehouse"
$ws4.Cells.Item(31,5).Value = "2024.03.09 14:00-03.09 15:30"
$ws4.Cells.Item(31,6).Value = 1
$ws4.Cells.Item(31,7).Value = "'380"
$ws4.Cells.Item(31,8).Value = $false
$ws4.Cells.Item(31,9).Value = "https://show.bilibili.com/platform/detail.html?id=81279"
$ws4.Cells.Item(31,10).Value = "//i2.hdslb.com/bfs/openplatform/202401/tMZ1Jp2G1705992352054.jpeg"

# 广州·春卷饭 十周年 2024 专场演出 (now on row 36 after the insert above) : 248 -> 300
$ws4.Cells.Item(36,6).Value = 300

# Renumber column A (id) sequentially 0..N for both changed sheets, matching
# the "index" convention used throughout the workbook.
for ($r = 1; $r -le $ws2.UsedRange.Rows.Count; $r++) {
    $ws2.Cells.Item($r,1).Value = $r - 1
}
for ($r = 1; $r -le $ws4.UsedRange.Rows.Count; $r++) {
    $ws4.Cells.Item($r,1).Value = $r - 1
}
